# Applies the "tuning run-through" edit described in the commit:
#  - flips the tune_this_time flag to TRUE on the workflow sheet
#  - highlights the add_trend row on the workflow sheet
#  - renames/retypes the "week" row on variables into a time_id/time pair (bold)
#  - adds a new "trend" row (13) on variables
#  - appends a new "Sheet1" worksheet with the same trend row content
#  - leaves the selection/active-sheet state the way the author's last save did

$wb = $excel.ActiveWorkbook

# ---- workflow sheet (sheet 1) ----
$wsWorkflow = $wb.Worksheets.Item("workflow")
$wsWorkflow.Activate()
$wsWorkflow.Range("B2").Value = $true
$wsWorkflow.Range("B2").Select()

# ---- variables sheet (sheet 2) ----
$wsVariables = $wb.Worksheets.Item("variables")
$wsVariables.Activate()
$wsVariables.Range("C4").Value = "time_id"
$wsVariables.Range("D4").Value = "time"
$wsVariables.Range("C4:D4").Font.Bold = $true

$wsVariables.Range("A13").Value = "trend"
$wsVariables.Range("B13").Value = "trend"
$wsVariables.Range("C13").Value = "predictor"
$wsVariables.Range("D13").Value = "trend"

# ---- back to workflow sheet for the highlight fill ----
# (done after the bold-font style above so new style indices are allocated
#  in the same order the saved workbook uses)
$wsWorkflow.Range("A6:C6").Interior.Color = 65535

# ---- role controls sheet (sheet 3) ----
$wsRoleControls = $wb.Worksheets.Item("role controls")
$wsRoleControls.Activate()
$wsRoleControls.Range("E28").Select()

# ---- new Sheet1 appended at the end ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNew.Range("A1").Value = "trend"
$wsNew.Range("B1").Value = "trend"
$wsNew.Range("C1").Value = "predictor"
$wsNew.Range("D1").Value = "trend"
$wsNew.Range("A1:D1").Select()

# ---- finish back on the variables sheet, matching the saved selection ----
$wsVariables.Activate()
$wsVariables.Range("A13:D13").Select()
